# Seahawks 2021 Team Data - log Week 17 game data
# (also fixes a couple of tiebreak-related DEF/ST totals per the Simulate_Season.py fix)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's per-play yardage logs to the running lists
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value() + " 3 -1 1 20 15 13 23 1 6 4 2 5 5 2 37 7 3 3 6 13 8 0 3 1 2 9 1 6 17 7 0 -2 0 6 3 4 1 30"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value() + " -1 6 2 1 1 2 0 3 -3 3 26 31 5 0 -1 3 2 0 14 0 1"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value() + " 12 7 4 8 9 -5 58 11 28 6 7 13 8 1 13 11 12 17 1 15"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value() + " 6 12 15 13 16 1 6 20 42 7 16 11 5 6 6 8 9 11 23 4 4"

# ---------------------------------------------------------------------------
# OFF sheet: Week 17 offensive down/distance + play totals
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("B2").Value = 4
$offWs.Range("C2").Value = 190
$offWs.Range("E2").Value = 10
$offWs.Range("F2").Value = 52
$offWs.Range("G2").Value = 51
$offWs.Range("J2").Value = 24
$offWs.Range("L2").Value = 230
$offWs.Range("M2").Value = 155
$offWs.Range("O2").Value = 12
$offWs.Range("P2").Value = 4
$offWs.Range("Q2").Value = 468

$offWs.Range("C3").Value = 145
$offWs.Range("D3").Value = 4
$offWs.Range("E3").Value = 34
$offWs.Range("F3").Value = 87
$offWs.Range("G3").Value = 35
$offWs.Range("H3").Value = 22
$offWs.Range("I3").Value = 54
$offWs.Range("J3").Value = 48
$offWs.Range("N3").Value = 24

# ---------------------------------------------------------------------------
# DEF sheet: Week 17 defensive down/distance + play totals
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 214
$defWs.Range("D2").Value = 14
$defWs.Range("E2").Value = 7
$defWs.Range("F2").Value = 65
$defWs.Range("G2").Value = 77
$defWs.Range("I2").Value = 11
$defWs.Range("J2").Value = 35
$defWs.Range("L2").Value = 313
$defWs.Range("M2").Value = 204
$defWs.Range("O2").Value = 29
$defWs.Range("P2").Value = 12
$defWs.Range("Q2").Value = 563

$defWs.Range("B3").Value = 12
$defWs.Range("C3").Value = 225
$defWs.Range("E3").Value = 34
$defWs.Range("F3").Value = 122
$defWs.Range("G3").Value = 43
$defWs.Range("I3").Value = 54
$defWs.Range("J3").Value = 55

# ---------------------------------------------------------------------------
# ST sheet: Week 17 special-teams totals + per-kick logs
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 73
$stWs.Range("D2").Value = 80
$stWs.Range("F2").Value = 260
$stWs.Range("G2").Value = 234
$stWs.Range("J2").Value = 91
$stWs.Range("K2").Value = 85
$stWs.Range("L2").Value = 61
$stWs.Range("M2").Value = 52
$stWs.Range("N2").Value = 37
$stWs.Range("O2").Value = 23

$stWs.Range("B3").Value = 26

$stWs.Range("B4").Value = $stWs.Range("B4").Value() + " 57 65"
$stWs.Range("B5").Value = $stWs.Range("B5").Value() + " 47 27"
$stWs.Range("B6").Value = $stWs.Range("B6").Value() + " 21 27"
$stWs.Range("D3").Value = $stWs.Range("D3").Value() + " 45"
$stWs.Range("D4").Value = $stWs.Range("D4").Value() + " 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value() + " 3 0 30"

# ---------------------------------------------------------------------------
# TURNS sheet: Week 17 turnover totals
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("C2").Value = 7
$turnsWs.Range("E3").Value = 7

# ---------------------------------------------------------------------------
# PEN sheet: Week 17 penalty totals
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 14
$penWs.Range("D3").Value = 10
